$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap (prolificid/index, name) pairs for row-pairs (2,3), (6,7), (9,10)
# while leaving A (index) and H (re_rank) columns untouched, and set new
# "realeffort" (F) values for every data row (2-13).

# Row 2 <-> Row 3
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("D2").Value = "Katherine"

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("D3").Value = "Melissa"

# Row 6 <-> Row 7
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "5eeaa065c7acf61c4322f6d9"
$ws.Range("D6").Value = "Yonifredy"

$ws.Range("B7").Value = 11
$ws.Range("C7").Value = "5f5ea8227fa75676f56f9276"
$ws.Range("D7").Value = "Carlos"

# Row 9 <-> Row 10
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "5e0adc8f4cac6834756db412"
$ws.Range("D9").Value = "Mary"

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "5e706891c396cc64388ef760"
$ws.Range("D10").Value = "Maria"

# Updated "realeffort" values (column F) for all data rows
$ws.Range("F2").Value = 8.316648944792245
$ws.Range("F3").Value = 8.023344841524992
$ws.Range("F4").Value = 7.429443214079729
$ws.Range("F5").Value = 7.094459853851288
$ws.Range("F6").Value = 6.304187637973969
$ws.Range("F7").Value = 6.19077550199683
$ws.Range("F8").Value = 5.065205973220809
$ws.Range("F9").Value = 3.358111939047832
$ws.Range("F10").Value = 3.344821734808749
$ws.Range("F11").Value = 2.475630392065158
$ws.Range("F12").Value = 1.281029176420817
$ws.Range("F13").Value = 0.379747773547242
